$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Sema6d", "Plxna1", "ECs", 3, 1, 24.44575933333333, 73.337278, 0.4034052273345712, 0.4034052273345712, 3, 1, 5.635297666666666, 16.905893, 0.09323918274550408, 0.09323918274550407, 137.7591305310282, 1239.832174779254, 0.03761317371193971, 0.03761317371193969),
    @("ECs", "Sema6d", "Plxna1", "FAPs", 3, 1, 24.44575933333333, 73.337278, 0.4034052273345712, 0.4034052273345712, 3, 1, 18.53178066666667, 55.595342, 0.3066187779927862, 0.3066187779927861, 453.0234501954528, 4077.211051759076, 0.1236916178412283, 0.1236916178412283),
    @("ECs", "Sema6d", "Plxna1", "M2", 3, 1, 24.44575933333333, 73.337278, 0.4034052273345712, 0.4034052273345712, 3, 1, 14.77436233333333, 44.323087, 0.2444501694549869, 0.2444501694549868, 361.1705059041317, 3250.534553137186, 0.09861247618096344, 0.0986124761809634),
    @("ECs", "Sema6d", "Plxna1", "sCs", 3, 1, 24.44575933333333, 73.337278, 0.4034052273345712, 0.4034052273345712, 3, 1, 21.497717, 64.493151, 0.355691869806723, 0.3556918698067229, 525.5280159981086, 4729.752143982978, 0.1434879596004398, 0.1434879596004398),
    @("FAPs", "Sema6d", "Plxna1", "ECs", 3, 1, 14.28901333333333, 42.86704, 0.2357980618855278, 0.2357980618855278, 3, 1, 5.635297666666666, 16.905893, 0.09323918274550408, 0.09323918274550407, 80.52284349630223, 724.70559146672, 0.02198561858318041, 0.02198561858318041),
    @("FAPs", "Sema6d", "Plxna1", "FAPs", 3, 1, 14.28901333333333, 42.86704, 0.2357980618855278, 0.2357980618855278, 3, 1, 18.53178066666667, 55.595342, 0.3066187779927862, 0.3066187779927861, 264.8008610364089, 2383.20774932768, 0.07230011358840792, 0.0723001135884079),
    @("FAPs", "Sema6d", "Plxna1", "M2", 3, 1, 14.28901333333333, 42.86704, 0.2357980618855278, 0.2357980618855278, 3, 1, 14.77436233333333, 44.323087, 0.2444501694549869, 0.2444501694549868, 211.1110603724978, 1899.99954335248, 0.05764087618507476, 0.05764087618507475),
    @("FAPs", "Sema6d", "Plxna1", "sCs", 3, 1, 14.28901333333333, 42.86704, 0.2357980618855278, 0.2357980618855278, 3, 1, 21.497717, 64.493151, 0.355691869806723, 0.3556918698067229, 307.1811648492267, 2764.63048364304, 0.08387145352886477, 0.08387145352886476),
    @("M2", "Sema6d", "Plxna1", "ECs", 2, 0.6666666666666666, 0.480005, 1.440015, 0.007921068169999337, 0.007921068169999336, 3, 1, 5.635297666666666, 16.905893, 0.09323918274550408, 0.09323918274550407, 2.704971056488333, 24.344739508395, 0.0007385539226421638, 0.0007385539226421635),
    @("M2", "Sema6d", "Plxna1", "FAPs", 2, 0.6666666666666666, 0.480005, 1.440015, 0.007921068169999337, 0.007921068169999336, 3, 1, 18.53178066666667, 55.595342, 0.3066187779927862, 0.3066187779927861, 8.895347378903333, 80.05812641013, 0.002428748242682752, 0.002428748242682751),
    @("M2", "Sema6d", "Plxna1", "M2", 2, 0.6666666666666666, 0.480005, 1.440015, 0.007921068169999337, 0.007921068169999336, 3, 1, 14.77436233333333, 44.323087, 0.2444501694549869, 0.2444501694549868, 7.091767791811668, 63.825910126305, 0.001936306456420841, 0.00193630645642084),
    @("M2", "Sema6d", "Plxna1", "sCs", 2, 0.6666666666666666, 0.480005, 1.440015, 0.007921068169999337, 0.007921068169999336, 3, 1, 21.497717, 64.493151, 0.355691869806723, 0.3556918698067229, 10.319011648585, 92.87110483726499, 0.002817459548253581, 0.002817459548253581),
    @("sCs", "Sema6d", "Plxna1", "ECs", 3, 1, 21.38374133333333, 64.151224, 0.3528756426099016, 0.3528756426099016, 3, 1, 5.635297666666666, 16.905893, 0.09323918274550408, 0.09323918274550407, 120.5037476403369, 1084.533728763032, 0.03290183652774181, 0.0329018365277418),
    @("sCs", "Sema6d", "Plxna1", "FAPs", 3, 1, 21.38374133333333, 64.151224, 0.3528756426099016, 0.3528756426099016, 3, 1, 18.53178066666667, 55.595342, 0.3066187779927862, 0.3066187779927861, 396.2788042220675, 3566.509237998608, 0.1081982983204672, 0.1081982983204672),
    @("sCs", "Sema6d", "Plxna1", "M2", 3, 1, 21.38374133333333, 64.151224, 0.3528756426099016, 0.3528756426099016, 3, 1, 14.77436233333333, 44.323087, 0.2444501694549869, 0.2444501694549868, 315.9311425009431, 2843.380282508488, 0.08626051063252783, 0.08626051063252782),
    @("sCs", "Sema6d", "Plxna1", "sCs", 3, 1, 21.38374133333333, 64.151224, 0.3528756426099016, 0.3528756426099016, 3, 1, 21.497717, 64.493151, 0.355691869806723, 0.3556918698067229, 459.7016195852026, 4137.314576266824, 0.1255149971291648, 0.1255149971291648),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $rowData[$j]
    }
}
